# Insert a new data row at row 519 (pushing existing rows 519:555 down to 520:556)
# and populate it with a new "Betarraga" price record for
# "Región de Arica y Parinacota" dated 44746 (2022-07-04).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("519:519").Insert()

$ws.Range("A519").Value = 9
$ws.Range("B519").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C519").Value = "Metropolitana"
$ws.Range("D519").Value = 44746
$ws.Range("E519").Value = 13
$ws.Range("F519").Value = 100114014
$ws.Range("G519").Value = "Betarraga"
$ws.Range("H519").Value = "Sin especificar"
$ws.Range("I519").Value = "Primera"
$ws.Range("J519").Value = 5200
$ws.Range("K519").Value = 110
$ws.Range("L519").Value = 120
$ws.Range("M519").Value = 115
$ws.Range("N519").Value = "`$/unidad"
$ws.Range("O519").Value = "Región de Arica y Parinacota"
$ws.Range("P519").Value = 115
$ws.Range("Q519").Value = 1
$ws.Range("R519").Value = "Hortaliza"
